# Add a new row for "FUSI-ZON CREAM 15 GM" (price 48.00) as item #8 at row 14,
# pushing the existing items (old rows 14-20) down by one row (new rows 15-21),
# then fix up the trailing total / footer rows and update the generated timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: shift the totals row (21) and footer row (22) down to 22 / 23. ---
# Grab their current formatting first (PasteSpecial formats only - keeps the same
# style ids instead of minting new ones), then move the values/merges.

$ws.Range("A21:Q21").Copy()
$ws.Range("A22:Q22").PasteSpecial(-4122)
$ws.Range("A22:Q22").RowHeight = $ws.Range("A21:Q21").RowHeight

$ws.Range("A22:Q22").Copy()
$ws.Range("A23:Q23").PasteSpecial(-4122)
$ws.Range("A23:Q23").RowHeight = $ws.Range("A22:Q22").RowHeight

# totals row values (P21/Q21 -> P22/Q22)
$ws.Cells.Item(22, 16).Value2 = $ws.Cells.Item(21, 16).Value2

# footer row values (row22 -> row23)
$ws.Cells.Item(23, 1).Value2 = $ws.Cells.Item(22, 1).Value2
$ws.Cells.Item(23, 7).Value2 = $ws.Cells.Item(22, 7).Value2
$ws.Cells.Item(23, 11).Value2 = $ws.Cells.Item(22, 11).Value2

# remove the old merges that are now stale, and the old footer/total merges
$ws.Range("P21:Q21").UnMerge()
$ws.Range("A22:F22").UnMerge()
$ws.Range("G22:I22").UnMerge()
$ws.Range("K22:Q22").UnMerge()

# re-create merges at the shifted rows
$ws.Range("P22:Q22").Merge()
$ws.Range("A23:F23").Merge()
$ws.Range("G23:I23").Merge()
$ws.Range("K23:Q23").Merge()

# --- 2. Turn the (now vacated) row 21 into a regular data row, matching the
#        formatting of the other data rows (copy format from row 20). ---
$ws.Range("A20:Q20").Copy()
$ws.Range("A21:Q21").PasteSpecial(-4122)
$ws.Range("A21:Q21").RowHeight = $ws.Range("A20:Q20").RowHeight

$ws.Range("A20:B20,C20:G20,H20:K20,L20:M20,N20:O20").Copy() | Out-Null

$ws.Range("A21:B21").Merge()
$ws.Range("C21:G21").Merge()
$ws.Range("H21:K21").Merge()
$ws.Range("L21:M21").Merge()
$ws.Range("N21:O21").Merge()

# --- 3. Shift the data rows 14-20 down to 15-21 (bottom-up so nothing is lost). ---
for ($r = 20; $r -ge 14; $r--) {
    $dst = $r + 1
    $ws.Cells.Item($dst, 1).Value2  = $r
    $ws.Cells.Item($dst, 3).Value2  = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($dst, 8).Value2  = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($dst, 12).Value2 = $ws.Cells.Item($r, 12).Value2
    $ws.Cells.Item($dst, 14).Value2 = $ws.Cells.Item($r, 14).Value2
    $ws.Cells.Item($dst, 16).Value2 = $ws.Cells.Item($r, 16).Value2
    $ws.Cells.Item($dst, 17).Value2 = $ws.Cells.Item($r, 17).Value2
}

# item # 8 -> 8 stays the same, the shift above already renumbered 9..15 into
# rows 15..21; now it's +1 too many on row 15 (it copied "8" from row14 via loop
# start), so fix row 15's number explicitly afterwards. (A21 should be 15.)
$ws.Cells.Item(21, 1).Value2 = 15

# --- 4. Write the new item into row 14. ---
$ws.Cells.Item(14, 1).Value2  = 8
$ws.Cells.Item(14, 3).Value2  = "FUSI-ZON CREAM 15 GM"
$ws.Cells.Item(14, 8).Value2  = "1:0"
$ws.Cells.Item(14, 12).Value2 = "1"
$ws.Cells.Item(14, 14).Value2 = "48.00"
$ws.Cells.Item(14, 16).Value2 = "96.0000"
$ws.Cells.Item(14, 17).Value2 = "2:0"

# --- 5. Update the grand total and the generated timestamp. ---
$ws.Cells.Item(22, 16).Value2 = 573.35

$ws.Cells.Item(23, 1).Value2 = "Tuesday, 5 August, 2025 10:57 AM"
